# Weekly fruit/vegetable price update:
# Insert a brand-new observation as row 493 (pushing the existing
# rows 493:503 down to 494:504, which is exactly what the diff shows -
# every old row's data reappears one row lower, unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row above the current row 493; Excel shifts 493:503 -> 494:504.
$ws.Rows.Item(493).Insert()

# Populate the newly inserted row 493 with the new record's data.
$ws.Cells.Item(493, 1).Value = 9
$ws.Cells.Item(493, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(493, 3).Value = "Metropolitana"
$ws.Cells.Item(493, 4).Value = 44939
$ws.Cells.Item(493, 5).Value = 13
$ws.Cells.Item(493, 6).Value = 100112032
$ws.Cells.Item(493, 7).Value = "Zapallo italiano"
$ws.Cells.Item(493, 8).Value = "Sin especificar"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 520
$ws.Cells.Item(493, 11).Value = 4000
$ws.Cells.Item(493, 12).Value = 5000
$ws.Cells.Item(493, 13).Value = 4500
$ws.Cells.Item(493, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(493, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(493, 16).Value = 90
$ws.Cells.Item(493, 17).Value = 50
$ws.Cells.Item(493, 18).Value = "Hortaliza"
